# Update the dSF column (column F) with refreshed data pulled from the
# source, and recompute the mean for the header/summary rows affected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 3
    5  = 8
    6  = 1
    7  = -1
    8  = 3
    9  = 2
    10 = 1
    11 = -2
    12 = 2
    13 = 2
    14 = 4
    15 = 5
    17 = 5
    18 = 1
    19 = -3
    20 = -1
    21 = 3
    22 = -2
    23 = 6
    24 = 7
    25 = 3
    26 = 2
    28 = -1
    29 = 1
    32 = 3
    33 = 10
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
